$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gasto de gasolina")

# --- Header row (row 1) ---
# A1 previously held an empty, styled cell; it disappears entirely.
$ws.Range("A1").Clear()

# B1 becomes "Lunes" and picks up the bold/centered-wrap style already used
# by C1:F1 (style index 3) instead of its former style (index 1).
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "Lunes"

# C1:F1 become the weekday headers (note trailing spaces / no accent on
# "Miercoles"), keeping their existing style.
$ws.Range("C1").Value = "Martes "
$ws.Range("D1").Value = "Miercoles "
$ws.Range("E1").Value = "Jueves "
$ws.Range("F1").Value = "Viernes "

# New blank, formatted cells H1:L1 (same style as the other header cells).
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("L1").PasteSpecial(-4122)

# --- Row labels (column A, rows 2:6) lose their bold style and get new text ---
$ws.Range("A2").Value = "Oficina "
$ws.Range("A2").Style = "Normal"

$ws.Range("A3").Value = "Naucalpan"
$ws.Range("A3").Style = "Normal"

$ws.Range("A4").Value = "Santa Fe"
$ws.Range("A4").Style = "Normal"

$ws.Range("A5").Value = "Lomas "
$ws.Range("A5").Style = "Normal"

$ws.Range("A6").Value = "Santa Fe"
$ws.Range("A6").Style = "Normal"

# --- Selection / dimension bookkeeping ---
$ws.Range("I19").Select()
